$wb = $excel.ActiveWorkbook

# Mapping of row -> new value for column F ("想去人数") that changed
$updates = @{
    3  = 6308
    4  = 177
    6  = 39
    7  = 1899
    8  = 1434
    10 = 960
    11 = 262
    12 = 5587
}

# Both "展览" and "全部类型" sheets contain identical data and both were updated
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
